# Validate Feedback 3 Business Validation :
# 3. Create a separate class with application business level validation
#
# Adds the "Business Validation" design notes (new rows on the "Design"
# sheet describing the validate/ValidationInterceptor/TaskRequestValidator/
# AppConfig pieces) and records the Feedback-3 write-up plus a missing
# bullet of explanatory text on the "2.Validation" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Design")
$ws2 = $wb.Worksheets.Item("Summary")
$ws4 = $wb.Worksheets.Item("2.Validation")

# --- "Design" sheet: three new rows (16-18) describing the interceptor-based
# business validation solution ---
$ws1.Range("C16").Value = "validate"
$ws1.Range("E16").Value = "HandlerInterceptor preHandle"
$ws1.Range("C17").Value = "validation"
$ws1.Range("D16").Value = "ValidationInterceptor"
$ws1.Range("D17").Value = "TaskRequestValidator"
$ws1.Range("E17").Value = "Business level valiadation which can be placed in controller or service"

# The "class name" cells (column D) in this table use the same font style as
# the existing TaskService/CreateTaskRequest rows (D6/D10) - copy that
# formatting across instead of leaving the default style.
$ws1.Range("D6").Copy()
$ws1.Range("D16").PasteSpecial(-4122)
$ws1.Range("D17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- "2.Validation" sheet: add the Feedback 3 write-up and the explanatory
# text that was missing next to the "3." bullet ---
$ws4.Range("D6").Value = "Validate Feedback 3 Business Validation : "
$ws4.Range("D21").Value = "Create a separate class with application business level validation"
$ws4.Range("E6").Value = "3. Create a separate class with application business level validation"
$ws4.Columns("D:D").ColumnWidth = 36.3
$ws4.Range("D11").Select()

# --- back on "Design": finish the table with the AppConfig/interceptor
# registration row, then widen column D to fit the new class names ---
$ws1.Range("C18").Value = "config"
$ws1.Range("D18").Value = "AppConfig"
$ws1.Range("E18").Value = "register interceptor"
$ws1.Columns("D:D").ColumnWidth = 23.15
$ws1.PageSetup.Orientation = 1

# --- "Summary" sheet: cursor moved here at some point too ---
$ws2.Range("C4").Select()

# --- finally land back on "Design", which becomes the active tab ---
$ws1.Select()
$ws1.Range("C19").Select()
